# Update Name of Algo
# Apply the numeric corrections described in the diff to Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("E3").Value = 12.5733
$ws.Range("B9").Value = 8.580000000000005
$ws.Range("E11").Value = 13.4093
$ws.Range("B18").Value = 4.593500000000003
$ws.Range("B20").Value = 5.717999999999997
